$d = $word.ActiveDocument

# 1) Remove the stale _GoBack bookmark that currently sits after "Nagercoil"
#    (Word moves this bookmark to the most recent edit location; it will be
#    re-created below at the new edit location.)
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2) Append " , Collaboration" right after "Good Communication"
$rng = $d.Content
$found = $rng.Find.Execute("Good Communication", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)
$rng.InsertAfter(" , Collaboration")
$rng.Font.Size = 11.5

# 3) Re-create the _GoBack bookmark at the new edit location (end of the
#    just-inserted text). Adding a zero-length bookmark exactly at a
#    paragraph boundary is unreliable, so temporarily pad with a marker
#    character, drop the bookmark right before it, then remove the marker -
#    the bookmark (being zero-width) stays put as the marker is deleted.
$padRng = $d.Range($rng.End, $rng.End)
$padRng.InsertAfter("X")

$bmRng = $d.Range($rng.End, $rng.End)
$d.Bookmarks.Add("_GoBack", $bmRng)

$markerRng = $d.Range($rng.End, $rng.End + 1)
$markerRng.Delete()
